$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$values = @(
    "7-6=1",
    "90-46=44",
    "35+18=53",
    "72+9=81",
    "17+69=86",
    "7+72=79",
    "20+77=97",
    "39+58=97",
    "15+82=97",
    "18+33=51",
    "56-55=1",
    "68-43=25",
    "26+20=46",
    "30+32=62",
    "91-49=42",
    "1+35=36",
    "74+23=97",
    "78-6=72",
    "7+83=90",
    "77-32=45",
    "51+48=99",
    "23-21=2",
    "6+56=62",
    "43-1=42",
    "69-0=69",
    "21+43=64",
    "65-13=52",
    "42+3=45",
    "51-2=49",
    "61+20=81",
    "56+11=67",
    "99-2=97",
    "38-25=13",
    "45-2=43",
    "30+16=46",
    "2+14=16",
    "37-8=29",
    "7+25=32",
    "34+44=78",
    "10+75=85",
    "0+27=27",
    "28+43=71",
    "57+26=83",
    "86-78=8",
    "41-3=38",
    "40+50=90",
    "43-5=38",
    "68-65=3",
    "38+61=99",
    "60+9=69",
    "36+18=54",
    "10+88=98",
    "50+4=54",
    "86-37=49",
    "38+23=61",
    "24+73=97",
    "98-96=2",
    "86-13=73",
    "78-44=34",
    "5+47=52",
    "87-0=87",
    "49+41=90",
    "56-20=36",
    "36+18=54",
    "57-24=33",
    "15-7=8",
    "35-16=19",
    "82+0=82",
    "51-46=5",
    "2+35=37",
    "68+10=78",
    "64-27=37",
    "30-3=27",
    "79-25=54",
    "75-4=71",
    "72-11=61",
    "51-42=9",
    "68-62=6",
    "91-51=40",
    "14+46=60",
    "57+10=67",
    "6+77=83",
    "92-39=53",
    "46+4=50",
    "9+72=81",
    "7+41=48",
    "79-14=65",
    "72+9=81",
    "59-23=36",
    "46+19=65",
    "26+28=54",
    "4+11=15",
    "28+4=32",
    "68-21=47",
    "84-61=23",
    "74+13=87",
    "41+33=74",
    "95-17=78",
    "5+57=62",
    "70-26=44"
)

$cols = 5
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = [int][Math]::Floor($i / $cols) + 1
    $col = ($i % $cols) + 1
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $values[$i]
}

Write-Output "done"